# Auto-generated edit script for 'Add data for 2022-10-26'
# Updates 2022 (column I) and, in a couple of places, 2021 (column H)
# crime-count figures across the Citywide Totals, By Neighborhood, and
# individual neighborhood sheets, per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H2").Value = 7242
$ws.Range("I2").Value = 6001
$ws.Range("I3").Value = 6251
$ws.Range("I4").Value = 1439
$ws.Range("I5").Value = 583
$ws.Range("I6").Value = 7088
$ws.Range("H7").Value = 25986
$ws.Range("I7").Value = 21362

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 61
$ws.Range("I3").Value = 57
$ws.Range("I4").Value = 36
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 228

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I5").Value = 17
$ws.Range("I6").Value = 197
$ws.Range("I7").Value = 681

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I3").Value = 141
$ws.Range("I6").Value = 97
$ws.Range("I7").Value = 383

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 192
$ws.Range("I3").Value = 304
$ws.Range("I7").Value = 825

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 166
$ws.Range("I6").Value = 142
$ws.Range("I7").Value = 501

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 166
$ws.Range("I6").Value = 155
$ws.Range("I7").Value = 668
$ws.Range("I8").Value = 1281
$ws.Range("I16").Value = 62
$ws.Range("I19").Value = 592
$ws.Range("I23").Value = 213
$ws.Range("I29").Value = 1319
$ws.Range("I33").Value = 966
$ws.Range("I37").Value = 681
$ws.Range("I42").Value = 738
$ws.Range("I43").Value = 181
$ws.Range("I46").Value = 45
$ws.Range("I47").Value = 149
$ws.Range("I52").Value = 465
$ws.Range("I54").Value = 436
$ws.Range("I58").Value = 12
$ws.Range("I60").Value = 116
$ws.Range("H63").Value = 225
$ws.Range("I63").Value = 79
$ws.Range("I65").Value = 501
$ws.Range("I67").Value = 825
$ws.Range("I76").Value = 307
$ws.Range("I77").Value = 137
$ws.Range("I78").Value = 291
$ws.Range("I79").Value = 604
$ws.Range("I83").Value = 462
$ws.Range("I85").Value = 974
$ws.Range("I89").Value = 248
$ws.Range("I90").Value = 262
$ws.Range("I91").Value = 229
$ws.Range("I94").Value = 222
$ws.Range("I95").Value = 326
$ws.Range("I96").Value = 228
$ws.Range("I97").Value = 182
$ws.Range("I99").Value = 383
$ws.Range("H101").Value = 25986
$ws.Range("I101").Value = 21362

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 158
$ws.Range("I3").Value = 167
$ws.Range("I4").Value = 19
$ws.Range("I7").Value = 462

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 326

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 215
$ws.Range("I3").Value = 367
$ws.Range("I6").Value = 304
$ws.Range("I7").Value = 966

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 95
$ws.Range("I7").Value = 436

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 384
$ws.Range("I7").Value = 1319

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 200
$ws.Range("I5").Value = 15
$ws.Range("I7").Value = 592

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I3").Value = 70
$ws.Range("I6").Value = 143
$ws.Range("I7").Value = 307

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 277
$ws.Range("I3").Value = 375
$ws.Range("I6").Value = 244
$ws.Range("I7").Value = 974

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 55
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 183
$ws.Range("I3").Value = 235
$ws.Range("I6").Value = 243
$ws.Range("I7").Value = 738

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I5").Value = 6
$ws.Range("I7").Value = 291

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 74
$ws.Range("I7").Value = 213

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I6").Value = 178
$ws.Range("I7").Value = 604

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 122
$ws.Range("I7").Value = 465

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 38
$ws.Range("I6").Value = 126
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I5").Value = 9
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 59
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I2").Value = 30
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 362
$ws.Range("I4").Value = 79
$ws.Range("I6").Value = 415
$ws.Range("I7").Value = 1281

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 262

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I6").Value = 102
$ws.Range("I7").Value = 181

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 219
$ws.Range("I3").Value = 207
$ws.Range("I7").Value = 668

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 12

